$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "293.25"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.92%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "30.62"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "4.06%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.159"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.58%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07133"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "7.10%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.535"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.99%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.613"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "6.12%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.398"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.81%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9171"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.15%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1641"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "4.11%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07904"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "21.91%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07781"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.85%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02948"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.37%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09001"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.24%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001577"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.13%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0006578"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.69%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006233"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.38%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.483"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.09%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.244"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.26%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3252"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.11%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "5.06%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.154"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.04%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.1590"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.49%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04539"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.79%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.56%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004240"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "2.42%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001168"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-6.41%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001688"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "4.44%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04421"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "4.80%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007034"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "4.64%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "2.11%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002207"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "12.17%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01338"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "8.28%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005851"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "5.31%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.897"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-3.51%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.01298"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.55%"
